$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset shrinks from 8 data rows (rows 2-9) to 5 data rows (rows 2-6).
# Remove the now-unused trailing rows first so the final dimension becomes A1:T6.
$ws.Rows("7:9").Delete()

# --- Row 2: ECs -> ECs (Efnb3/Rhbdl2) ---
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb3"
$ws.Range("C2").Value = "Rhbdl2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.3050205
$ws.Range("H2").Value = 0.6100410000000001
$ws.Range("I2").Value = 0.2027672883226855
$ws.Range("J2").Value = 0.1761044107725932
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.183143
$ws.Range("N2").Value = 6.366286000000001
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.9709238694315002
$ws.Range("R2").Value = 3.883695477726001
$ws.Range("S2").Value = 0.2027672883226855
$ws.Range("T2").Value = 0.1761044107725932

# --- Row 3: FAPs -> ECs (Efnb3/Rhbdl2) ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efnb3"
$ws.Range("C3").Value = "Rhbdl2"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2207143333333333
$ws.Range("H3").Value = 0.662143
$ws.Range("I3").Value = 0.1467234066692216
$ws.Range("J3").Value = 0.1911450260920121
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.183143
$ws.Range("N3").Value = 6.366286000000001
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.7025652851496668
$ws.Range("R3").Value = 4.215391710898
$ws.Range("S3").Value = 0.1467234066692216
$ws.Range("T3").Value = 0.1911450260920121

# --- Row 4: Inflammatory-Mac -> ECs (Efnb3/Rhbdl2) ---
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Efnb3"
$ws.Range("C4").Value = "Rhbdl2"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06966766666666667
$ws.Range("H4").Value = 0.209003
$ws.Range("I4").Value = 0.04631270309296832
$ws.Range("J4").Value = 0.06033422370742998
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.183143
$ws.Range("N4").Value = 6.366286000000001
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.2217621454763334
$ws.Range("R4").Value = 1.330572872858
$ws.Range("S4").Value = 0.04631270309296832
$ws.Range("T4").Value = 0.06033422370742998

# --- Row 5: MuSCs -> ECs (Efnb3/Rhbdl2) ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Efnb3"
$ws.Range("C5").Value = "Rhbdl2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7437579999999999
$ws.Range("H5").Value = 1.487516
$ws.Range("I5").Value = 0.4944251052906407
$ws.Range("J5").Value = 0.4294106932071856
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.183143
$ws.Range("N5").Value = 6.366286000000001
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 2.367488071394
$ws.Range("R5").Value = 9.469952285575999
$ws.Range("S5").Value = 0.4944251052906407
$ws.Range("T5").Value = 0.4294106932071856

# --- Row 6: Neutrophils -> ECs (Efnb3/Rhbdl2) ---
$ws.Range("A6").Value = "Neutrophils"
$ws.Range("B6").Value = "Efnb3"
$ws.Range("C6").Value = "Rhbdl2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.165128
$ws.Range("H6").Value = 0.495384
$ws.Range("I6").Value = 0.1097714966244839
$ws.Range("J6").Value = 0.1430056462207791
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.183143
$ws.Range("N6").Value = 6.366286000000001
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 0.525626037304
$ws.Range("R6").Value = 3.153756223824
$ws.Range("S6").Value = 0.1097714966244839
$ws.Range("T6").Value = 0.1430056462207791
